# Upload of data after editor comments
# - Rename promoter labels in the header row (C1/D1), swapping which
#   label lands in which column so the final text is "P_relB" / "P_bolA".
# - Move the active selection from G6 to C2.
# - Reposition the workbook window (xWindow/yWindow) to match the
#   author's saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the promoter columns.
$ws.Range("C1").Value = "P_relB"
$ws.Range("D1").Value = "P_bolA"

# Update the saved selection/active cell.
$ws.Range("C2").Select()

# Update the saved window position.
$win = $wb.Windows.Item(1)
$win.Left = 28800
$win.Top = -1800
